# Weekly update: insert a new price observation as the first data row
# (row 268), pushing the existing rows 268-382 down to 269-383.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(268).Insert()

$ws.Range("A268").Value = 10
$ws.Range("B268").Value = "Vega Modelo de Temuco"
$ws.Range("C268").Value = "La Araucanía"
$ws.Range("D268").Value = 44825
$ws.Range("E268").Value = 9
$ws.Range("F268").Value = 100112044
$ws.Range("G268").Value = "Perejil"
$ws.Range("H268").Value = "Sin especificar"
$ws.Range("I268").Value = "Primera"
$ws.Range("J268").Value = 50
$ws.Range("K268").Value = 3300
$ws.Range("L268").Value = 3300
$ws.Range("M268").Value = 3300
$ws.Range("N268").Value = "$/docena de atados (3 kilos)"
$ws.Range("O268").Value = "Región Metropolitana"
$ws.Range("P268").Value = 1100
$ws.Range("Q268").Value = 3
$ws.Range("R268").Value = "Hortaliza"
